$d = $word.ActiveDocument

# Update the ID placeholder text in the first paragraph, absorbing the
# trailing space run so it disappears entirely.
$d.Content.Find.Execute("**ID__AFFARS_5332_topic_14__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5332_7__ID**", 2)

# Adjust the first paragraph's formatting: indent + a thin paragraph border
# on all four sides (5pt gap).
$p = $d.Paragraphs(1)
$p.Range.ParagraphFormat.LeftIndent = 11.25
$p.Borders.DistanceFromTop = 5
$p.Borders.DistanceFromLeft = 5
$p.Borders.DistanceFromBottom = 5
$p.Borders.DistanceFromRight = 5
